# The sheet's monthly rows (2014-01 .. 2017-12) get re-ordered: within each
# calendar year's 12-row block, October/November/December move to the front
# of the block, ahead of January..September (labels travel together with
# their whole row of data). Row 1 (header) and the overall layout otherwise
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 49
$numCols      = 7   # columns A..G

# Snapshot every data row (A..G) before any writes, so later writes in the
# same pass never clobber a value we still need to read.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Each block is 12 rows (one calendar year). Within a block, old position
# 0..11 = Jan..Dec; new position order is Oct,Nov,Dec,Jan,Feb,...,Sep.
$blockStarts = @(2, 14, 26, 38)
$rotate = @(9, 10, 11, 0, 1, 2, 3, 4, 5, 6, 7, 8)

$mapping = @{}
foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt 12; $i++) {
        $newRow = $start + $i
        $oldRow = $start + $rotate[$i]
        $mapping[$newRow] = $oldRow
    }
}

for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
    $oldRow = $mapping[$newRow]
    $srcRow = $snapshot[$oldRow]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $srcRow[$c - 1]
    }
}
